# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" column (E16:E71) lists 56 monthly periods
# (1612 = Dec-2016 ... 2107 = Jul-2021). The periods are reordered to
# newest-first (reverse chronological order) — i.e. the old periods are
# "removed" from the top and the newest ones are added/moved to the top of
# the account-statement list. The underlying "database" (the F column
# "Salario Basico" amounts) is also touched: the two rows that used to be
# out of step with the rest of the column (120000) swap places with the one
# outlier (100000).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Chronological order of the periods exactly as they appeared (oldest first).
$periods = @("1612","1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712","1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812","1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912","2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102","2103","2104","2105","2106","2107")

# Reverse to newest-first order.
$periods = $periods[($periods.Length - 1)..0]

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value2 = $periods[$i]
}

# Swap the two "Salario Basico" values that differ from the common 120000:
# F16 (120000 -> 100000) and F71 (100000 -> 120000).
$ws.Range("F16").Value2 = 100000
$ws.Range("F71").Value2 = 120000
